# Append the latest gold-price row (row 84) to Sheet1, mirroring the
# existing rows: column A holds a dd-mm-yyyy date *label* stored as text,
# column B holds the descriptive price sentence (also text).
#
# A plain `.Value = "09-12-2025"` assignment gets auto-recognised as a
# date literal by the engine (turning the cell into a numeric date serial
# and forcing a brand-new, inconsistent style). To keep the cell as plain
# text - matching every other date-label cell in the column - we compute
# the string via a throwaway formula cell (formula results are never
# re-interpreted as dates) and bring only its *value* across with
# PasteSpecial(xlPasteValues), which copies the underlying text verbatim
# without re-running Excel's "smart" type detection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDateLabel = "09-12-2025"
$newPriceText = "The price of gold in India today is ₹13,009 per gram for 24 karat gold, ₹11,925 per gram for 22 karat gold and ₹9,757 per gram for 18 karat gold (also called 999 gold)."

$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""$newDateLabel"""
$scratch.Copy()
$ws.Range("A84").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B84").Value = $newPriceText
